$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132..223 down to 133..224
$ws.Rows("132:132").Insert()

# Populate the newly inserted row 132 with the new record's data
$ws.Range("A132").Value = 10
$ws.Range("B132").Value = "Vega Modelo de Temuco"
$ws.Range("C132").Value = "La Araucanía"
$ws.Range("D132").Value = 45068
$ws.Range("E132").Value = 9
$ws.Range("F132").Value = "Fruta"
$ws.Range("G132").Value = 100104
$ws.Range("H132").Value = "Frutos de pepita"
$ws.Range("I132").Value = 100104001
$ws.Range("J132").Value = "Granada"
$ws.Range("K132").Value = "Wonderfull"
$ws.Range("L132").Value = "Primera"
$ws.Range("M132").Value = 110
$ws.Range("N132").Value = 22000
$ws.Range("O132").Value = 22000
$ws.Range("P132").Value = 22000
$ws.Range("Q132").Value = "$/bandeja 15 kilos granel"
$ws.Range("R132").Value = "Región de O'Higgins"
$ws.Range("S132").Value = 1467
$ws.Range("T132").Value = 15
